$wb = $excel.ActiveWorkbook

# --- locate the existing sheets ---
$summary = $wb.Worksheets.Item(1)   # "总计"
$q3 = $wb.Worksheets.Item(2)        # currently "2022-Q3" (the original 72-fund data)

# Helper: write a value into a Range as a literal text string (never auto-coerced to a
# number by Excel), without touching / re-minting the cell's existing style. We do this by
# writing a `="<text>"` formula into a scratch cell far outside any used range, copying it,
# and pasting *values only* into the destination (PasteSpecial xlPasteValues = -4163).
function Set-TextValue {
    param($sheet, $range, [string]$text)
    $scratch = $sheet.Range("ZZ1000")
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Clear()
}

# --- Step 1: duplicate the original "2022-Q3" sheet to a new sheet right after it ---
# This new sheet keeps the original Q3 fund-holding data + formatting untouched.
$q3copy = $wb.Worksheets.Add($null, $q3)
$q3.UsedRange.Copy($q3copy.Range("A1"))
# A1 was never populated in the source data (the used range just starts there); drop the
# stray empty cell the range-copy leaves behind so the sheet matches the source exactly.
$q3copy.Range("A1").Clear()

# --- Step 2: the original sheet becomes the new "2022-Q4" sheet; the copy keeps the "2022-Q3" name ---
$q3.Name = "2022-Q4"
$q3copy.Name = "2022-Q3"

# --- Step 3: clear the old Q3 content out of what is now the Q4 sheet ---
$q3.Cells.Clear()

# --- Step 4: write the new Q4 header row (text, non-numeric-looking so no coercion risk) ---
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Apply the bold/centered "header + row-number column" style (style used by B1:D1/A2 on the
# "总计" sheet) onto the new header row and the A-column of the data rows, by copying the
# format (not the value) from the summary sheet's already-styled cells.
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q3.Range("A2:A15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 5: new Q4 fund-holding data rows ---
# Columns B, D, E, F, G hold numeric-looking values that must stay literal TEXT (as in the
# source data) rather than being coerced to numbers; column C (fund name) is never
# numeric-looking so a plain assignment is safe; column A / H are genuine numbers.
$q4data = @(
    @(0, "180031", "银华中小盘精选混合",       "37.32", "91.61", "4.75", "1.7727", 6),
    @(1, "005543", "银华心诚灵活配置混合A",     "18.70", "89.69", "3.13", "0.5853", 9),
    @(2, "519001", "银华核心价值优选混合",       "24.89", "90.13", "1.90", "0.4729", 10),
    @(3, "009085", "银华丰享一年持有期混合",     "3.60",  "91.96", "4.73", "0.1703", 6),
    @(4, "001534", "华宝万物互联灵活配置混合A",  "0.77",  "91.08", "4.37", "0.0336", 2),
    @(5, "002861", "工银智能制造股票",           "0.61",  "93.74", "4.99", "0.0304", 10),
    @(6, "004258", "国寿安保稳嘉混合A",          "2.16",  "23.32", "1.20", "0.0259", 7),
    @(7, "011543", "中加科瑞混合A",              "0.98",  "26.11", "1.37", "0.0134", 9),
    @(8, "014042", "银华心诚灵活配置混合C",      "0.26",  "89.69", "3.13", "0.0081", 9),
    @(9, "014014", "招商臻选平衡混合A",          "0.25",  "66.99", "2.38", "0.0060", 10),
    @(10, "014015", "招商臻选平衡混合C",         "0.19",  "66.99", "2.38", "0.0045", 10),
    @(11, "004259", "国寿安保稳嘉混合C",         "0.01",  "23.32", "1.20", "0.0001", 7),
    @(12, "016463", "华宝万物互联灵活配置混合C", "0.00",  "91.08", "4.37", $null, 2),
    @(13, "011544", "中加科瑞混合C",             "0.00",  "26.11", "1.37", $null, 9)
)

$r = 2
foreach ($row in $q4data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $q3 $q3.Cells.Item($r, 2) $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    Set-TextValue $q3 $q3.Cells.Item($r, 4) $row[3]
    Set-TextValue $q3 $q3.Cells.Item($r, 5) $row[4]
    Set-TextValue $q3 $q3.Cells.Item($r, 6) $row[5]
    if ($row[6] -eq $null) {
        $q3.Cells.Item($r, 7).Value = 0
    } else {
        Set-TextValue $q3 $q3.Cells.Item($r, 7) $row[6]
    }
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Match the "总计" sheet's page margins (0.75in / 1in / 0.5in) on the new Q4 sheet.
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# --- Step 6: update the "总计" (summary) sheet ---
# Row 2 becomes the new 2022-Q4 summary entry.
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 14
$summary.Range("D2").Value = 3.12

# Row 3 is added with the original 2022-Q3 summary entry (shifted down from row 2).
# Copy A2's formatting (bold/centered style) down to A3 first, then set its value.
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 71
$summary.Range("D3").Value = 27.78
